$wb = $excel.ActiveWorkbook

# --- Sheet: "Sampling Events" ---
# The survey actually covered three villages (DESA YENDE, DESA MENARBU,
# DESA KAYOP) under a single sampling event (PR001) rather than three
# separate sampling events (PR001/PR002/PR003). Consolidate the locality
# text on the surviving row, then drop the two redundant event rows.
$ws1 = $wb.Worksheets.Item("Sampling Events")
$ws1.Range("S2").Value = "DESA YENDE, DESA MENARBU, DESA KAYOP"
$ws1.Range("A3:A4").EntireRow.Delete()

# --- Sheet: "Occurrences" ---
# Renumber the occurrence/eventID identifiers: what used to be vouchers
# "VE001".."VE011" under PR001 plus "VE001" under a (now removed) PR002
# become a single consecutive "EM001".."EM012" sequence under PR001.
$ws2 = $wb.Worksheets.Item("Occurrences")
for ($i = 1; $i -le 12; $i++) {
    $row = $i + 1
    $num = "{0:D3}" -f $i
    $ws2.Cells.Item($row, 2).Value = "UNIPA -2006ES-AF002-PR001-EM$num"
}
